$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2033333333333333
$ws.Range("C2").Value = 0.5366666666666666
$ws.Range("J2").Value = 0.01666666666666667
$ws.Range("P2").Value = 0.14
$ws.Range("S2").Value = 0.1033333333333333
$ws.Range("B3").Value = 0.006134969325153374
$ws.Range("C3").Value = 0.01226993865030675
$ws.Range("J3").Value = 0.0245398773006135
$ws.Range("P3").Value = 0.8098159509202454
$ws.Range("S3").Value = 0.147239263803681
$ws.Range("J4").Value = 0.06060606060606061
$ws.Range("P4").Value = 0.7575757575757576
$ws.Range("S4").Value = 0.1818181818181818
$ws.Range("B6").Value = 0.04854368932038835
$ws.Range("D6").Value = 0.009708737864077669
$ws.Range("F6").Value = 0.05339805825242718
$ws.Range("J6").Value = 0.2233009708737864
$ws.Range("O6").Value = 0.01941747572815534
$ws.Range("Q6").Value = 0.2427184466019418
$ws.Range("R6").Value = 0.06796116504854369
$ws.Range("S6").Value = 0.3349514563106796
$ws.Range("B7").Value = 0.07051282051282051
$ws.Range("D7").Value = 0.03205128205128205
$ws.Range("F7").Value = 0.03846153846153846
$ws.Range("J7").Value = 0.1282051282051282
$ws.Range("O7").Value = 0.03205128205128205
$ws.Range("Q7").Value = 0.2115384615384615
$ws.Range("R7").Value = 0.08333333333333333
$ws.Range("S7").Value = 0.4038461538461539
$ws.Range("B8").Value = 0.1
$ws.Range("D8").Value = 0.02727272727272727
$ws.Range("F8").Value = 0.06818181818181818
$ws.Range("J8").Value = 0.1045454545454545
$ws.Range("O8").Value = 0.01363636363636364
$ws.Range("Q8").Value = 0.2568181818181818
$ws.Range("R8").Value = 0.09545454545454546
$ws.Range("S8").Value = 0.3340909090909091
$ws.Range("B9").Value = 0.1173708920187793
$ws.Range("D9").Value = 0.0187793427230047
$ws.Range("F9").Value = 0.06572769953051644
$ws.Range("J9").Value = 0.136150234741784
$ws.Range("O9").Value = 0.02347417840375587
$ws.Range("Q9").Value = 0.2394366197183098
$ws.Range("R9").Value = 0.07042253521126761
$ws.Range("S9").Value = 0.3286384976525822
$ws.Range("B10").Value = 0.1068100358422939
$ws.Range("D10").Value = 0.03010752688172043
$ws.Range("F10").Value = 0.06666666666666667
$ws.Range("J10").Value = 0.1197132616487455
$ws.Range("O10").Value = 0.02078853046594982
$ws.Range("Q10").Value = 0.2874551971326165
$ws.Range("R10").Value = 0.07670250896057347
$ws.Range("S10").Value = 0.2917562724014337
$ws.Range("G11").Value = 0.1422764227642276
$ws.Range("J11").Value = 0.07723577235772358
$ws.Range("K11").Value = 0.1951219512195122
$ws.Range("L11").Value = 0.573170731707317
$ws.Range("S11").Value = 0.01219512195121951
$ws.Range("G12").Value = 0.723404255319149
$ws.Range("J12").Value = 0.2340425531914894
$ws.Range("K12").Value = 0.02836879432624113
$ws.Range("S12").Value = 0.01418439716312057
$ws.Range("G13").Value = 0.7575757575757576
$ws.Range("J13").Value = 0.2424242424242424
$ws.Range("F15").Value = 0.004048582995951417
$ws.Range("H15").Value = 0.1214574898785425
$ws.Range("I15").Value = 0.05668016194331984
$ws.Range("J15").Value = 0.4210526315789473
$ws.Range("K15").Value = 0.03238866396761134
$ws.Range("M15").Value = 0.008097165991902834
$ws.Range("O15").Value = 0.07692307692307693
$ws.Range("S15").Value = 0.2793522267206478
$ws.Range("F16").Value = 0.02272727272727273
$ws.Range("H16").Value = 0.1863636363636364
$ws.Range("I16").Value = 0.06818181818181818
$ws.Range("J16").Value = 0.4318181818181818
$ws.Range("K16").Value = 0.08181818181818182
$ws.Range("M16").Value = 0.01363636363636364
$ws.Range("O16").Value = 0.06363636363636363
$ws.Range("S16").Value = 0.1318181818181818
$ws.Range("F17").Value = 0.009302325581395349
$ws.Range("H17").Value = 0.1426356589147287
$ws.Range("I17").Value = 0.09767441860465116
$ws.Range("J17").Value = 0.4666666666666667
$ws.Range("K17").Value = 0.09302325581395349
$ws.Range("M17").Value = 0.01395348837209302
$ws.Range("N17").Value = 0.001550387596899225
$ws.Range("O17").Value = 0.05891472868217054
$ws.Range("S17").Value = 0.1162790697674419
$ws.Range("F18").Value = 0.01036269430051814
$ws.Range("H18").Value = 0.2124352331606218
$ws.Range("I18").Value = 0.08808290155440414
$ws.Range("J18").Value = 0.4248704663212435
$ws.Range("K18").Value = 0.06217616580310881
$ws.Range("M18").Value = 0.01036269430051814
$ws.Range("O18").Value = 0.07772020725388601
$ws.Range("S18").Value = 0.1139896373056995
$ws.Range("F19").Value = 0.008042895442359249
$ws.Range("H19").Value = 0.2073279714030384
$ws.Range("I19").Value = 0.09472743521000894
$ws.Range("J19").Value = 0.4003574620196604
$ws.Range("K19").Value = 0.08668453976764968
$ws.Range("M19").Value = 0.01697944593386953
$ws.Range("N19").Value = 0.0008936550491510277
$ws.Range("O19").Value = 0.08132260947274352
$ws.Range("S19").Value = 0.1036639857015192
